# Append two newly-scraped Lancers listings to the top of the "ランサーズ"
# sheet's data table (pushing all existing listings down by two rows), and
# log the 2025-08-31 06:23:08 scrape run's stats as a new row at the bottom
# of the "統計" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "ランサーズ" - insert 2 fresh rows right below the header row.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ランサーズ")

$ws1.Rows.Item(2).Resize(2).Insert()

# New row 2: 【急募】小型BLE音声モジュールの試作開発依頼
$ws1.Range("A2").Value = "2025-08-31 06:23:08"
$ws1.Range("B2").Value = "【急募】小型BLE音声モジュールの試作開発依頼"
$ws1.Range("C2").Value = "システム開発"
$ws1.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws1.Range("E2").Value = "期限情報なし"
$ws1.Range("F2").Value = "https://www.lancers.jp/work/detail/5383341"
$ws1.Range("F2").Style = "Hyperlink"
$ws1.Range("G2").Value = 68
$ws1.Range("H2").Value = "◆開発"

# New row 3: LINE公式アカウントの自動応答・ステップ配信設定(文章提供あり)
$ws1.Range("A3").Value = "2025-08-31 06:23:08"
$ws1.Range("B3").Value = "LINE公式アカウントの自動応答・ステップ配信設定(文章提供あり)"
$ws1.Range("C3").Value = "システム開発"
$ws1.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws1.Range("E3").Value = "期限情報なし"
$ws1.Range("F3").Value = "https://www.lancers.jp/work/detail/5383334"
$ws1.Range("F3").Style = "Hyperlink"
$ws1.Range("G3").Value = 18
# (this listing has no matched skill keywords, so H3 is intentionally left blank)

# The row-insert above shifted every old row down by two, but it left the
# hyperlink relationships pinned to their original cell addresses, so the
# two rows that slid into the previously-unused F71/F72 slots ended up
# without a clickable link. Re-attach them.
$ws1.Hyperlinks.Add($ws1.Range("F71"), "https://www.lancers.jp/work/detail/5380357")
$ws1.Hyperlinks.Add($ws1.Range("F72"), "https://www.lancers.jp/work/detail/5380420")

# ---------------------------------------------------------------------
# Sheet 2: "統計" - append this run's summary stats.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("統計")

$nextRow = $ws2.UsedRange.Rows.Count + 1

$ws2.Cells.Item($nextRow, 1).Value = "2025-08-31T06:23:08.108535"
$ws2.Cells.Item($nextRow, 2).Value = 7
$ws2.Cells.Item($nextRow, 3).Value = "全案件リスト"
$ws2.Cells.Item($nextRow, 4).Value = 71.40000000000001
$ws2.Cells.Item($nextRow, 5).Value = 2
$ws2.Cells.Item($nextRow, 6).Value = 3
$ws2.Cells.Item($nextRow, 7).Value = 7
